# Trade #78 closed at 2026-02-17 21:13:40 - unknown UNKNOWN +0.000%
#
# This records a MarketMaking trade (internal MarketMaking-sheet trade #106,
# shown as row 107 on "All Trades" / row 74 on "MarketMaking") being closed
# via an early exit, and a brand-new MarketMaking trade (#139) being opened
# right after (row 140 on "All Trades" / row 107 on "MarketMaking").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 106      # Total Trades
$summary.Range("B9").Value = 46.23    # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row = row 5)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D5").Value = 73        # Trades
$status.Range("G5").Value = 47.95     # Win Rate %

# ---------------------------------------------------------------------------
# All Trades sheet - close out existing trade #106 (row 107)
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("G107").Value = 0.02        # Exit Price
$allTrades.Range("H107").Value = "CLOSED"    # Status
$allTrades.Range("K107").Value = 101.14      # Capital After
$allTrades.Range("L107").Value = "early_exit" # Exit Reason
$allTrades.Range("M107").Value = 0.13        # Duration (min)

# All Trades sheet - append newly opened trade #139 (row 140)
$allTrades.Range("A140").Value = 139
$allTrades.Range("B140").NumberFormat = "@"   # keep date-looking text as text
$allTrades.Range("B140").Value = "2026-02-17"
$allTrades.Range("C140").Value = "21:13:33"
$allTrades.Range("D140").Value = "MarketMaking"
$allTrades.Range("E140").Value = "DOWN"
$allTrades.Range("F140").Value = 0.02
$allTrades.Range("H140").Value = "OPEN"
$allTrades.Range("I140").Value = 0
$allTrades.Range("J140").Value = 0
$allTrades.Range("K140").Value = 101.1396151053151
$allTrades.Range("M140").Value = 0
$allTrades.Range("N140").Value = 0
$allTrades.Range("O140").Value = 0
$allTrades.Range("P140").Value = 0.6
$allTrades.Range("Q140").Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------------
# MarketMaking sheet - close out existing trade #106 (row 74)
# ---------------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Range("G74").Value = 0.02          # Exit Price
$mm.Range("H74").Value = "CLOSED"      # Status
$mm.Range("K74").Value = 101.14        # Capital After
$mm.Range("P74").Value = "early_exit"  # Exit Reason
$mm.Range("Q74").Value = 0.13          # Duration (min)

# MarketMaking sheet - append newly opened trade #139 (row 107)
$mm.Range("A107").Value = 139
$mm.Range("B107").NumberFormat = "@"          # keep date-looking text as text
$mm.Range("B107").Value = "2026-02-17"
$mm.Range("C107").Value = "21:13:33"
$mm.Range("D107").Value = "MarketMaking"
$mm.Range("E107").Value = "DOWN"
$mm.Range("F107").Value = 0.02
$mm.Range("H107").Value = "OPEN"
$mm.Range("I107").Value = 0
$mm.Range("J107").Value = 0
$mm.Range("K107").Value = 101.1396151053151
$mm.Range("L107").Value = 0
$mm.Range("M107").Value = 0
$mm.Range("N107").Value = 0.6
$mm.Range("O107").Value = "Normal spread capture: 19600 bps"
$mm.Range("Q107").Value = 0
